# Applies the Wed Aug  7 15:56:49 UTC 2024 cryptos list update (GitHub Actions)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "55.997.01"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.00%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.401.73"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -3.96%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "479.92"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.08"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +1.93%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.24%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.501"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -1.68%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.407.85"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -4.41%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0981"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.78%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.46"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -3.41%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.323"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -2.39%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.12%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.813.07"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -3.95%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "56.387.69"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.69%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.39"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -3.08%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000133"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.83%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.392.87"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -4.89%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.50"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.93%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "317.19"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.90%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.81"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -4.40%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.998"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.02%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.70"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.83%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "56.89"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -2.63%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.37%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.396"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -3.31%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.159"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -4.22%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.497.16"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -4.08%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.34"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -2.76%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0777"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.97%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.08%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "148.95"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.50%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.01"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.73%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.02%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.03"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -3.71%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -3.29%  "
$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = "Fetch.AI"
$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.848"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -2.19%  "
$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = "NEARProtocol"
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.60"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -2.97%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "33.55"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -2.00%  "
$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.36"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +2.93%  "
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "FirstDigitalUSD"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.999"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.52%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.40"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -4.12%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0544"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.74%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0948"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +4.79%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -5.26%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.50%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.67"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -2.95%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "255.53"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -2.54%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.58%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "17.06"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -3.44%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.781.08"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -7.55%  "

Write-Output "Applied cryptos update"
